$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Martes 21/05/2024"
$ws.Range("B3").Value = "Miércoles 22/05/2024"
$ws.Range("B4").Value = "Jueves 23/05/2024"
$ws.Range("B5").Value = "Lunes 27/05/2024"
$ws.Range("B6").Value = "Martes 28/05/2024"
$ws.Range("B7").Value = "Miércoles 29/05/2024"
$ws.Range("B8").Value = "Jueves 30/05/2024"
$ws.Range("B9").Value = "Viernes 31/05/2024"
$ws.Range("B10").Value = "Lunes 03/06/2024"
$ws.Range("B11").Value = "Martes 04/06/2024"
